# Update "Horarios" workbook for Linea 141 with the 05:57:08 scrape.
#
# Sheet 1 (LP1912) gains 5 new schedule rows (total filas 18 -> 23), and the
# "Ultima actualizacion" / "Total filas" headers on all three sheets move to
# the new scrape timestamp.

$wb = $excel.ActiveWorkbook

$newStamp = "05:57:08"

# ---------------------------------------------------------------------
# Sheet 1: LP1912 - new rows inserted + header counters updated
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: $newStamp"
$ws1.Range("A3").Value = "Total filas: 23"

# Insert the 5 new rows top-to-bottom so each insertion point is still
# expressed in terms of the *current* (already-shifted) row numbers.
$ws1.Rows.Item(17).Insert()
$ws1.Range("A17").Value = $newStamp
$ws1.Range("B17").Value = "07:01"
$ws1.Range("C17").Value = "16_SANTA ANA"
$ws1.Range("D17").Value = 64
$ws1.Range("E17").Value = "LP1912"

$ws1.Rows.Item(22).Insert()
$ws1.Range("A22").Value = $newStamp
$ws1.Range("B22").Value = "07:29"
$ws1.Range("C22").Value = "14_ABASTO"
$ws1.Range("D22").Value = 92
$ws1.Range("E22").Value = "LP1912"

$ws1.Rows.Item(24).Insert()
$ws1.Range("A24").Value = $newStamp
$ws1.Range("B24").Value = "07:34"
$ws1.Range("C24").Value = "23_HERNANDEZ"
$ws1.Range("D24").Value = 97
$ws1.Range("E24").Value = "LP1912"

$ws1.Rows.Item(27).Insert()
$ws1.Range("A27").Value = $newStamp
$ws1.Range("B27").Value = "07:44"
$ws1.Range("C27").Value = "10_OLMOS"
$ws1.Range("D27").Value = 107
$ws1.Range("E27").Value = "LP1912"

$ws1.Rows.Item(28).Insert()
$ws1.Range("A28").Value = $newStamp
$ws1.Range("B28").Value = "07:51"
$ws1.Range("C28").Value = "15_ABASTO"
$ws1.Range("D28").Value = 114
$ws1.Range("E28").Value = "LP1912"

# ---------------------------------------------------------------------
# Sheet 2 (LP1912-215) and Sheet 3 (6203-6173): only the timestamp banner
# changes, row data is identical.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: $newStamp"

$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: $newStamp"
